$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13.54093077557314
$ws.Range("E2").Value = 9.024927476408708
$ws.Range("F2").Value = 15.24477239439056
$ws.Range("G2").Value = 10.53971510077907
